# Rewrite 0.06: II_Windows, WPF...
# - WPF: Reimplemented start-up args[] input
# - WPF: Reimplemented menuitem shortcut keys
#
# The localization key list is reshuffled: the menu-related keys
# (File/LoadSimulation/SaveSimulation/ExitProgram/Help/AboutProgram and
# DeviceOptions/PauseDevice/NumericRowAmounts/TracingRowAmounts/FontSize/
# ColorScheme/ToggleFullscreen/CloseDevice/PatientOptions/NewPatient/
# EditPatient) are renamed with a "Menu" prefix: the "MenuFile" group now
# sits right under the header row, and the rest of the renamed menu keys
# move to the very end of the key list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that must end up blank (previously had content, now cleared).
$blankRows = @(2, 9, 10, 32, 42, 43)
foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 1).ClearContents()
}

# Cells whose final text already existed somewhere in the workbook before
# this edit (vital-sign / device labels that simply moved rows). Order of
# assignment amongst these does not matter.
$retained = @{
    1  = "Key"
    11 = "HeartRate"
    12 = "BloodPressure"
    13 = "RespiratoryRate"
    14 = "PulseOximetry"
    15 = "Temperature"
    16 = "EndTidalCO2"
    17 = "ArterialBloodPressure"
    18 = "CentralVenousPressure"
    19 = "PulmonaryArteryPressure"
    20 = "RespiratoryRhythm"
    21 = "InspiratoryExpiratoryRatio"
    22 = "CardiacRhythm"
    23 = "VitalSigns"
    24 = "AdvancedHemodynamics"
    25 = "RespiratoryProfile"
    26 = "CardiacProfile"
    27 = "UseDefaultVitalSignRanges"
    28 = "STSegmentElevation"
    29 = "TWaveElevation"
    30 = "ApplyChanges"
    31 = "ResetParameters"
    33 = "Devices"
    34 = "CardiacMonitor"
    35 = "12LeadECG"
    36 = "Defibrillator"
    37 = "Ventilator"
    38 = "IABP"
    39 = "Cardiotocograph"
    40 = "IVPump"
    41 = "LabResults"
}
foreach ($r in $retained.Keys) {
    $ws.Cells.Item($r, 1).Value = $retained[$r]
}

# Brand-new "Menu*" labels. These must be written in this exact order so
# that the workbook's underlying shared-string table is built in the same
# sequence as the authored edit (DeviceOptions group first, then the File
# group), matching the row layout: row 44 up through row 54, then rows 3-8.
$newLabelsInOrder = @(
    @{ Row = 44; Text = "MenuDeviceOptions" },
    @{ Row = 46; Text = "MenuNumericRowAmounts" },
    @{ Row = 45; Text = "MenuPauseDevice" },
    @{ Row = 47; Text = "MenuTracingRowAmounts" },
    @{ Row = 48; Text = "MenuFontSize" },
    @{ Row = 49; Text = "MenuColorScheme" },
    @{ Row = 50; Text = "MenuToggleFullscreen" },
    @{ Row = 51; Text = "MenuCloseDevice" },
    @{ Row = 52; Text = "MenuPatientOptions" },
    @{ Row = 53; Text = "MenuNewPatient" },
    @{ Row = 54; Text = "MenuEditPatient" },
    @{ Row = 3;  Text = "MenuFile" },
    @{ Row = 4;  Text = "MenuLoadSimulation" },
    @{ Row = 5;  Text = "MenuSaveSimulation" },
    @{ Row = 6;  Text = "MenuExitProgram" },
    @{ Row = 7;  Text = "MenuHelp" },
    @{ Row = 8;  Text = "MenuAboutProgram" }
)
foreach ($entry in $newLabelsInOrder) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Text
}

# Restore the frozen-pane view, but move the selection/scroll position to
# row 9 (whole row selected), matching the post-edit view state.
$ws.Range("A9:XFD9").Select()
